$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen columns A and B (values chosen so the persisted column width,
# after the host's internal character/pixel rounding, lands on
# A=37 and B as close as possible to 74.140625)
$ws.Columns.Item(1).ColumnWidth = 36.1666666666667
$ws.Columns.Item(2).ColumnWidth = 73.3333333333333

# Move selection to H26
$ws.Range("H26").Select()
